# add.feature of search in chatbot

$wb = $excel.ActiveWorkbook

# --- Assets sheet: append new row 11 ---
$assets = $wb.Worksheets.Item("Assets")
$assets.Range("A11").Value = "A1001"
$assets.Range("B11").Value = "Nil"
$assets.Range("C11").Value = "nil"

# D11/E11 look numeric/date-like ("12-1-2026", "7500") but must stay plain
# text, matching the inlineStr cells in the target file. Force text entry
# with a leading apostrophe, then strip the resulting "Text" number format
# so the cell keeps the default style.
$assets.Range("D11").Value = "'12-1-2026"
$assets.Range("D11").ClearFormats()
$assets.Range("E11").Value = "'7500"
$assets.Range("E11").ClearFormats()

$assets.Range("F11").Value = "IT"
$assets.Range("G11").Value = "Sanz"
$assets.Range("H11").Value = "HP"
$assets.Range("I11").Value = "Available"
$assets.Range("J11").Value = "Blooms"

# --- Transactions sheet: update row 10 values ---
$transactions = $wb.Worksheets.Item("Transactions")
$transactions.Range("B10").Value = "E007"
$transactions.Range("D10").Value = "tfg"
$transactions.Range("E10").Value = "tguj"
$transactions.Range("F10").Value = "ygkj"
$transactions.Range("G10").Value = "hfgh,m"
$transactions.Range("H10").Value = "jgfjhj,"
$transactions.Range("I10").Value = "hhjk"
$transactions.Range("J10").Value = "vbnm"
